# Add a new "2021" column (R) to the tourism-share-in-GDP table, mirroring
# the formatting already used by the adjacent "2020" column (Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (year header): copy Q4's formatting onto R4, then set the new year.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R4").Value = 2021

# Row 5 (data value): copy Q5's formatting onto R5, then set the new value.
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R5").Value = 3.6

$excel.CutCopyMode = $false

# Move the active selection, matching the saved view state.
$ws.Range("O9").Select()
